# Swap the species-record data between row 6 and row 7, keeping the
# shared/common columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AT, AW, AX, AY) untouched since they hold identical values in both rows.
#
# Note: this runtime's Range.Value getter is unreliable, so Value2 is used
# for both reads and writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry row-specific data which needs to be swapped between
# row 6 and row 7.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr6 = "$col`6"
    $addr7 = "$col`7"
    $v6 = $ws.Range($addr6).Value2
    $v7 = $ws.Range($addr7).Value2
    $ws.Range($addr6).Value2 = $v7
    $ws.Range($addr7).Value2 = $v6
}

# Column M ("Aktivitet") only had a value on row 7 ("äldre spår"); after the
# edit it belongs to row 6 instead, and row 7's M cell becomes empty.
$ws.Range("M6").Value2 = "äldre spår"
$ws.Range("M7").Value2 = ""
